$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.72352147102356
$ws.Range("B1").Value = 3.458930015563965
$ws.Range("C1").Value = 3.972855091094971
$ws.Range("D1").Value = 3.458053350448608
$ws.Range("E1").Value = 1.109869360923767
